# Generate Report for Handoff
#
# Adds a new handoff record for fd62fa24-6f8c-4c6f-a5e5-f1bd7b51c2d0.md as a
# new row at the bottom of each of the three tables (Overview, zh-cn, de-de),
# mirroring the existing cbb7c5ad-... row / record.

$wb = $excel.ActiveWorkbook

$newBase = "fd62fa24-6f8c-4c6f-a5e5-f1bd7b51c2d0"
$newMd = "$newBase.md"
$newDisplay = "e2e\$newBase.md"
$newUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3deb80433a69465009f720b105cdae1b4104f6a/e2e/$newMd"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | Path And Name | Extension | Publish URL |
#                    zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rOverview = $rowOverview.Range.Row

$wsOverview.Range("A$rOverview").Value = $newMd
$wsOverview.Range("C$rOverview").Value = ".md"
$wsOverview.Range("D$rOverview").Value = ""
$wsOverview.Range("E$rOverview").Value = "Ready for handoff"
$wsOverview.Range("F$rOverview").Value = "Ready for handoff"
$wsOverview.Range("G$rOverview").NumberFormat = $dateFmt
$wsOverview.Range("G$rOverview").Value = "2016-08-23 00:37:57"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B$rOverview"), $newUrl, "", "", $newDisplay)

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Source Path |
#   Priority | Content Duplicate | Latest Handoff File | Latest Handoff Datetime |
#   Latest Target File | Latest Handback File | Latest Handback DateTime |
#   Reference Tokens | To be localized | Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rZhCn = $rowZhCn.Range.Row

$wsZhCn.Range("B$rZhCn").Value = ".md"
$wsZhCn.Range("C$rZhCn").Value = "Ready for handoff"
$wsZhCn.Range("D$rZhCn").Value = "e2e"
$wsZhCn.Range("E$rZhCn").Value = "ht"
$wsZhCn.Range("F$rZhCn").Value = "'False"
$wsZhCn.Range("G$rZhCn").Value = "$newBase.2bd9fd2b0fdaee16b328d3058fd29ce5e1e4b71d.zh-cn.xlf"
$wsZhCn.Range("H$rZhCn").NumberFormat = $dateFmt
$wsZhCn.Range("H$rZhCn").Value = "2016-08-23 00:37:52"
$wsZhCn.Range("I$rZhCn").Value = ""
$wsZhCn.Range("J$rZhCn").Value = ""
$wsZhCn.Range("K$rZhCn").NumberFormat = $dateFmt
$wsZhCn.Range("K$rZhCn").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L$rZhCn").Value = ""
$wsZhCn.Range("M$rZhCn").Value = "'True"
$wsZhCn.Range("N$rZhCn").Value = ""
$wsZhCn.Range("O$rZhCn").Value = "'False"
$wsZhCn.Range("P$rZhCn").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A$rZhCn"), $newUrl, "", "", $newMd)

# ---------------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rDeDe = $rowDeDe.Range.Row

$wsDeDe.Range("B$rDeDe").Value = ".md"
$wsDeDe.Range("C$rDeDe").Value = "Ready for handoff"
$wsDeDe.Range("D$rDeDe").Value = "e2e"
$wsDeDe.Range("E$rDeDe").Value = "ht"
$wsDeDe.Range("F$rDeDe").Value = "'False"
$wsDeDe.Range("G$rDeDe").Value = "$newBase.2bd9fd2b0fdaee16b328d3058fd29ce5e1e4b71d.de-de.xlf"
$wsDeDe.Range("H$rDeDe").NumberFormat = $dateFmt
$wsDeDe.Range("H$rDeDe").Value = "2016-08-23 00:37:57"
$wsDeDe.Range("I$rDeDe").Value = ""
$wsDeDe.Range("J$rDeDe").Value = ""
$wsDeDe.Range("K$rDeDe").NumberFormat = $dateFmt
$wsDeDe.Range("K$rDeDe").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L$rDeDe").Value = ""
$wsDeDe.Range("M$rDeDe").Value = "'True"
$wsDeDe.Range("N$rDeDe").Value = ""
$wsDeDe.Range("O$rDeDe").Value = "'False"
$wsDeDe.Range("P$rDeDe").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A$rDeDe"), $newUrl, "", "", $newMd)

Write-Output "Handoff row added for $newMd"
